$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.885
$ws.Range("D3").Value = -7.58
$ws.Range("A4").Value = -20.862
$ws.Range("B4").Value = 7.043000000000001
$ws.Range("D4").Value = -7.388999999999998
$ws.Range("B5").Value = 6.619999999999999
$ws.Range("A6").Value = -20.97
$ws.Range("A7").Value = -21.038
$ws.Range("B8").Value = 6.347
$ws.Range("D9").Value = -7.631
$ws.Range("D11").Value = -8.513999999999999
$ws.Range("D14").Value = -7.694999999999999
$ws.Range("A16").Value = -21.119
$ws.Range("B16").Value = 6.611
$ws.Range("D18").Value = -8.210000000000003
$ws.Range("A20").Value = -22.232
$ws.Range("E20").Value = 12.922
$ws.Range("B22").Value = 7.202000000000001
$ws.Range("D25").Value = -8.412000000000001
